# Generate Report for handback
# Updates the zh-cn and de-de localization-status sheets to reflect a
# completed handback: status text changes from "Ready for handoff" to
# "Handed back: in sync with en-US", and a "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" are recorded for
# each of the two source files in each language sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 : 09cf1117-62f3-4130-aa59-b00a50001643.md
$wsZh.Range("B2").Value = $newStatus
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/576e41c7a88be9a45006703af2fea9efcc55ee98/e2e/09cf1117-62f3-4130-aa59-b00a50001643.md", "", "", "09cf1117-62f3-4130-aa59-b00a50001643.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ba1541c89540c4fcd2773a1bcb3925f41eb3c333/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/09cf1117-62f3-4130-aa59-b00a50001643.6948b8e253bdeb612cb6e16789f274eeb6ea7b25.zh-cn.xlf", "", "", "09cf1117-62f3-4130-aa59-b00a50001643.6948b8e253bdeb612cb6e16789f274eeb6ea7b25.zh-cn.xlf")
$wsZh.Range("G2").Value = "2016-01-26 05:37:07"

# Row 3 : 4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.md
$wsZh.Range("B3").Value = $newStatus
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/576e41c7a88be9a45006703af2fea9efcc55ee98/e2e/4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.md", "", "", "4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ba1541c89540c4fcd2773a1bcb3925f41eb3c333/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.2ae90447bf1701606f56496466a12aeea19087ff.zh-cn.xlf", "", "", "4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.2ae90447bf1701606f56496466a12aeea19087ff.zh-cn.xlf")
$wsZh.Range("G3").Value = "2016-01-26 05:37:07"

# Give the two new hyperlink columns (E: Latest Target File, F: Latest
# Handback File) the same visual style as the other hyperlink columns
# (A: Source File Name, C: Latest Handoff File) on this sheet - single
# underline, cornflower-blue text (FF6495ED), matching the workbook's
# custom "HyperLink" cell style.
$wsZh.Range("E2:F3").Font.Underline = 2
$wsZh.Range("E2:F3").Font.Color = 15570276

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 : 09cf1117-62f3-4130-aa59-b00a50001643.md
$wsDe.Range("B2").Value = $newStatus
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/576e41c7a88be9a45006703af2fea9efcc55ee98/e2e/09cf1117-62f3-4130-aa59-b00a50001643.md", "", "", "09cf1117-62f3-4130-aa59-b00a50001643.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/910acae7f7f5f02827399ba07ef85dd744087f58/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/09cf1117-62f3-4130-aa59-b00a50001643.6948b8e253bdeb612cb6e16789f274eeb6ea7b25.de-de.xlf", "", "", "09cf1117-62f3-4130-aa59-b00a50001643.6948b8e253bdeb612cb6e16789f274eeb6ea7b25.de-de.xlf")
$wsDe.Range("G2").Value = "2016-01-26 05:37:24"

# Row 3 : 4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.md
$wsDe.Range("B3").Value = $newStatus
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/576e41c7a88be9a45006703af2fea9efcc55ee98/e2e/4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.md", "", "", "4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/910acae7f7f5f02827399ba07ef85dd744087f58/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.2ae90447bf1701606f56496466a12aeea19087ff.de-de.xlf", "", "", "4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.2ae90447bf1701606f56496466a12aeea19087ff.de-de.xlf")
$wsDe.Range("G3").Value = "2016-01-26 05:37:24"

$wsDe.Range("E2:F3").Font.Underline = 2
$wsDe.Range("E2:F3").Font.Color = 15570276
